# Update the worksheet date header and every division-problem answer
# cell in the table to the values for the new day's worksheet.
# Find.Execute args: (FindText, MatchCase, MatchWholeWord, MatchWildcards,
#   MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format, ReplaceWith,
#   Replace). Wrap=1 (wdFindContinue), Replace=2 (wdReplaceAll).
$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-05-31 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-01 Sunday", 2) | Out-Null
$d.Content.Find.Execute("29÷5=5, 4", $true, $false, $false, $false, $false, $true, 1, $false, "48÷5=9, 3", 2) | Out-Null
$d.Content.Find.Execute("67÷2=33, 1", $true, $false, $false, $false, $false, $true, 1, $false, "96÷9=10, 6", 2) | Out-Null
$d.Content.Find.Execute("61÷3=20, 1", $true, $false, $false, $false, $false, $true, 1, $false, "49÷5=9, 4", 2) | Out-Null
$d.Content.Find.Execute("18÷7=2, 4", $true, $false, $false, $false, $false, $true, 1, $false, "52÷2=26, 0", 2) | Out-Null
$d.Content.Find.Execute("54÷6=9, 0", $true, $false, $false, $false, $false, $true, 1, $false, "29÷9=3, 2", 2) | Out-Null
$d.Content.Find.Execute("73÷8=9, 1", $true, $false, $false, $false, $false, $true, 1, $false, "27÷9=3, 0", 2) | Out-Null
$d.Content.Find.Execute("94÷2=47, 0", $true, $false, $false, $false, $false, $true, 1, $false, "77÷8=9, 5", 2) | Out-Null
$d.Content.Find.Execute("87÷9=9, 6", $true, $false, $false, $false, $false, $true, 1, $false, "70÷3=23, 1", 2) | Out-Null
$d.Content.Find.Execute("10÷8=1, 2", $true, $false, $false, $false, $false, $true, 1, $false, "49÷3=16, 1", 2) | Out-Null
$d.Content.Find.Execute("46÷7=6, 4", $true, $false, $false, $false, $false, $true, 1, $false, "78÷8=9, 6", 2) | Out-Null
$d.Content.Find.Execute("17÷4=4, 1", $true, $false, $false, $false, $false, $true, 1, $false, "16÷5=3, 1", 2) | Out-Null
$d.Content.Find.Execute("79÷8=9, 7", $true, $false, $false, $false, $false, $true, 1, $false, "25÷2=12, 1", 2) | Out-Null
$d.Content.Find.Execute("43÷6=7, 1", $true, $false, $false, $false, $false, $true, 1, $false, "38÷2=19, 0", 2) | Out-Null
$d.Content.Find.Execute("91÷5=18, 1", $true, $false, $false, $false, $false, $true, 1, $false, "77÷9=8, 5", 2) | Out-Null
$d.Content.Find.Execute("63÷7=9, 0", $true, $false, $false, $false, $false, $true, 1, $false, "69÷8=8, 5", 2) | Out-Null
$d.Content.Find.Execute("68÷4=17, 0", $true, $false, $false, $false, $false, $true, 1, $false, "24÷2=12, 0", 2) | Out-Null
$d.Content.Find.Execute("98÷4=24, 2", $true, $false, $false, $false, $false, $true, 1, $false, "41÷4=10, 1", 2) | Out-Null
$d.Content.Find.Execute("39÷3=13, 0", $true, $false, $false, $false, $false, $true, 1, $false, "11÷4=2, 3", 2) | Out-Null
$d.Content.Find.Execute("59÷5=11, 4", $true, $false, $false, $false, $false, $true, 1, $false, "14÷9=1, 5", 2) | Out-Null
$d.Content.Find.Execute("32÷4=8, 0", $true, $false, $false, $false, $false, $true, 1, $false, "92÷4=23, 0", 2) | Out-Null
$d.Content.Find.Execute("28÷9=3, 1", $true, $false, $false, $false, $false, $true, 1, $false, "77÷5=15, 2", 2) | Out-Null
$d.Content.Find.Execute("82÷4=20, 2", $true, $false, $false, $false, $false, $true, 1, $false, "28÷5=5, 3", 2) | Out-Null
$d.Content.Find.Execute("22÷3=7, 1", $true, $false, $false, $false, $false, $true, 1, $false, "17÷6=2, 5", 2) | Out-Null
$d.Content.Find.Execute("15÷6=2, 3", $true, $false, $false, $false, $false, $true, 1, $false, "66÷5=13, 1", 2) | Out-Null
$d.Content.Find.Execute("34÷4=8, 2", $true, $false, $false, $false, $false, $true, 1, $false, "43÷5=8, 3", 2) | Out-Null
